$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-07-06 Sunday"

# Update each arithmetic problem cell in the table, in document order (row-major)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "8+13="
$t.Cell(1,2).Range.Text = "6+5="
$t.Cell(1,3).Range.Text = "71-68="
$t.Cell(1,4).Range.Text = "70-46="
$t.Cell(1,5).Range.Text = "45+29="
$t.Cell(2,1).Range.Text = "71-35="
$t.Cell(2,2).Range.Text = "6+67="
$t.Cell(2,3).Range.Text = "40-31="
$t.Cell(2,4).Range.Text = "6+38="
$t.Cell(2,5).Range.Text = "54-27="
$t.Cell(3,1).Range.Text = "61-5="
$t.Cell(3,2).Range.Text = "53+19="
$t.Cell(3,3).Range.Text = "36+56="
$t.Cell(3,4).Range.Text = "19+39="
$t.Cell(3,5).Range.Text = "40-8="
$t.Cell(4,1).Range.Text = "56+25="
$t.Cell(4,2).Range.Text = "91-28="
$t.Cell(4,3).Range.Text = "17+28="
$t.Cell(4,4).Range.Text = "93-39="
$t.Cell(4,5).Range.Text = "48+36="
$t.Cell(5,1).Range.Text = "16+17="
$t.Cell(5,2).Range.Text = "35-17="
$t.Cell(5,3).Range.Text = "64-29="
$t.Cell(5,4).Range.Text = "21-16="
$t.Cell(5,5).Range.Text = "6+35="
$t.Cell(6,1).Range.Text = "59+25="
$t.Cell(6,2).Range.Text = "27+6="
$t.Cell(6,3).Range.Text = "53+29="
$t.Cell(6,4).Range.Text = "40-39="
$t.Cell(6,5).Range.Text = "62-23="
$t.Cell(7,1).Range.Text = "17+5="
$t.Cell(7,2).Range.Text = "72-54="
$t.Cell(7,3).Range.Text = "9+32="
$t.Cell(7,4).Range.Text = "90-64="
$t.Cell(7,5).Range.Text = "9+22="
$t.Cell(8,1).Range.Text = "61-49="
$t.Cell(8,2).Range.Text = "41-17="
$t.Cell(8,3).Range.Text = "46+28="
$t.Cell(8,4).Range.Text = "6+76="
$t.Cell(8,5).Range.Text = "42+29="
$t.Cell(9,1).Range.Text = "9+77="
$t.Cell(9,2).Range.Text = "5+27="
$t.Cell(9,3).Range.Text = "9+6="
$t.Cell(9,4).Range.Text = "27+66="
$t.Cell(9,5).Range.Text = "79+4="
$t.Cell(10,1).Range.Text = "32-7="
$t.Cell(10,2).Range.Text = "84-15="
$t.Cell(10,3).Range.Text = "6+69="
$t.Cell(10,4).Range.Text = "56-8="
$t.Cell(10,5).Range.Text = "14+28="
$t.Cell(11,1).Range.Text = "82-45="
$t.Cell(11,2).Range.Text = "87-68="
$t.Cell(11,3).Range.Text = "63-28="
$t.Cell(11,4).Range.Text = "84-68="
$t.Cell(11,5).Range.Text = "34-25="
$t.Cell(12,1).Range.Text = "49+43="
$t.Cell(12,2).Range.Text = "9+76="
$t.Cell(12,3).Range.Text = "57+15="
$t.Cell(12,4).Range.Text = "60-36="
$t.Cell(12,5).Range.Text = "21-19="
$t.Cell(13,1).Range.Text = "71-38="
$t.Cell(13,2).Range.Text = "90-15="
$t.Cell(13,3).Range.Text = "58-29="
$t.Cell(13,4).Range.Text = "45-9="
$t.Cell(13,5).Range.Text = "61-26="
$t.Cell(14,1).Range.Text = "59+7="
$t.Cell(14,2).Range.Text = "63-48="
$t.Cell(14,3).Range.Text = "70-19="
$t.Cell(14,4).Range.Text = "49+29="
$t.Cell(14,5).Range.Text = "31-24="
$t.Cell(15,1).Range.Text = "33-29="
$t.Cell(15,2).Range.Text = "67+24="
$t.Cell(15,3).Range.Text = "27+37="
$t.Cell(15,4).Range.Text = "93-68="
$t.Cell(15,5).Range.Text = "96-38="
$t.Cell(16,1).Range.Text = "76-49="
$t.Cell(16,2).Range.Text = "70-37="
$t.Cell(16,3).Range.Text = "73-55="
$t.Cell(16,4).Range.Text = "27+47="
$t.Cell(16,5).Range.Text = "64-28="
$t.Cell(17,1).Range.Text = "27-18="
$t.Cell(17,2).Range.Text = "29+39="
$t.Cell(17,3).Range.Text = "9+67="
$t.Cell(17,4).Range.Text = "84-28="
$t.Cell(17,5).Range.Text = "53+18="
$t.Cell(18,1).Range.Text = "66+7="
$t.Cell(18,2).Range.Text = "81-48="
$t.Cell(18,3).Range.Text = "64-9="
$t.Cell(18,4).Range.Text = "5+9="
$t.Cell(18,5).Range.Text = "58+4="
$t.Cell(19,1).Range.Text = "75-27="
$t.Cell(19,2).Range.Text = "62-3="
$t.Cell(19,3).Range.Text = "58+9="
$t.Cell(19,4).Range.Text = "19+64="
$t.Cell(19,5).Range.Text = "24+48="
$t.Cell(20,1).Range.Text = "72-37="
$t.Cell(20,2).Range.Text = "3+48="
$t.Cell(20,3).Range.Text = "89+7="
$t.Cell(20,4).Range.Text = "19+45="
$t.Cell(20,5).Range.Text = "32+59="
